$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new product ("شمع حريمي") was sold and needs to be added to the report,
# in the position currently occupied by row 25 ("كالونا "). The existing
# row 25 data is pushed down to a newly-inserted row 26, and the totals /
# footer rows shift down by one (26->27, 27->28). The grand total and the
# generated-at timestamp are also refreshed.
# ---------------------------------------------------------------------------

# 1) Insert a new blank row at 26 - this shifts the old row 26 (totals) to
#    27 and the old row 27 (footer) to 28, updating merged ranges below.
$ws.Rows.Item(26).Insert()

# 2) Clone row 25's cell formatting (styles/borders/fill/number formats)
#    into the new row 26 so it matches the rest of the product rows.
$ws.Range("A25:Q25").Copy()
$ws.Range("A26:Q26").PasteSpecial(-4122)

# 3) Clone row 25's values (and value *types*, important so numeric-looking
#    text like "15.0000" stays text) into row 26 - this duplicates row 25
#    ("كالونا ") down into row 26.
$ws.Range("A25:Q25").Copy()
$ws.Range("A26:Q26").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# 4) Row 26 is the 20th product line (row 25 stays the 19th).
$ws.Range("A26").Value = 20

# 5) Merge the row-26 cell groups the same way every other product row is
#    merged.
$ws.Range("A26:B26").Merge()
$ws.Range("C26:G26").Merge()
$ws.Range("H26:K26").Merge()
$ws.Range("L26:M26").Merge()
$ws.Range("N26:O26").Merge()

# 6) Row 26 (just like the neighbouring rows) autofits to 25.5pt.
$ws.Rows.Item(26).RowHeight = 25.5

# 7) Overwrite row 25 in place with the new product's data.
$ws.Range("C25").Value = "شمع حريمي"
$ws.Range("H25").Value = "8:0"
$ws.Range("N25").Value = "50.00"

# P25 keeps a numeric-style format (0.00) but its stored value is text in
# this workbook ("50.0000"), so route the text through a scratch cell
# (well outside the used range) formatted as Text and paste *values only*
# in order to keep P25's own style (s=11) while still storing a text value.
$helper = $ws.Range("A100")
$helper.NumberFormat = "@"
$helper.Value = "50.0000"
$helper.Copy()
$ws.Range("P25").PasteSpecial(-4163)
$helper.Clear()
$excel.CutCopyMode = 0

# 8) Grand total increases by the new product's price (866.495 + 50.00).
$ws.Range("P27").Value = 916.495

# 9) The footer timestamp (now row 28) is refreshed to the new export time.
$ws.Range("A28").Value = "Friday, 19 September, 2025 3:30 PM"
